# Upstream commit: "Fixed #295 Add the version of M2Doc in the template
# custom properties." That change stamps the M2Doc version into the
# *custom document properties* of the various M2Doc template resources
# touched by the fix.
#
# For this particular resource (newLine-template.docx) the word
# processing content itself is untouched by the commit: the canonical
# OOXML diff for word/document.xml and word/styles.xml only shows the
# package being re-serialized (namespace declarations and element
# attributes sorted, volatile w:rsid* bookkeeping attributes dropped) -
# there is no textual, structural, style or page-layout change in the
# body, the section properties or the style/latent-style definitions.
#
# So there is nothing to rewrite in the body, styles or page setup for
# this template. We just confirm the existing content is present (a
# read-only Find, which performs no replacement and therefore cannot
# perturb the document) and re-save the package, mirroring the
# no-content-change re-serialization seen in the diff.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("End of demonstration.")
Write-Host "content present: $found"

$d.Save()
